$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Add the new meeting entry captured on 02/21/2022 into row 3.
# ---------------------------------------------------------------

# Date (02/21/2022)
$ws.Range("A3").NumberFormat = "mm/dd/yyyy"
$ws.Range("A3").Value2 = 44613

# Start / End time
$ws.Range("B3").NumberFormat = "h:mm am/pm"
$ws.Range("B3").Value2 = 0.6923611111111111
$ws.Range("C3").NumberFormat = "h:mm am/pm"
$ws.Range("C3").Value2 = 0.6965277777777777

# Duration = ABS(start - end), same formula pattern used in row 2
$ws.Range("D3").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D3").Formula = "=ABS(B3-C3)"

# Client / Group + Attendees
$ws.Range("E3").Value2 = "Client"
$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("F3").Value2 = "Harrison, Nick, Jacob, Daniel, Brian"

# ---------------------------------------------------------------
# Update the totals row label from "Total Minutes:" to "Total Time:"
# The SUM formula in F20 recalculates automatically now that D3
# has a value.
# ---------------------------------------------------------------
$ws.Range("E20").Value2 = "Total Time:"
